# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# de-de handback file is now in sync with en-US (no more "stale handback"
# error), refreshes the handback timestamps, and widens a couple of columns
# that now show longer status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet ----------------------------------------------------
# The zh-cn / de-de status mirrors shown on the Overview tab
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# --- zh-cn sheet -----------------------------------------------------
# Status -> "Handed back: in sync with en-US"
$wsZhCn.Range("C2").Value = $newStatus
# Latest Handback DateTime refreshed
$wsZhCn.Range("K2").Value = "2016-08-24 20:51:30"
# Error Detail cleared (handback version now in sync, no error)
$wsZhCn.Range("P2").Value = ""

# --- de-de sheet -------------------------------------------------------
# Status -> "Handed back: in sync with en-US"
$wsDeDe.Range("C2").Value = $newStatus
# Latest Handback DateTime refreshed
$wsDeDe.Range("K2").Value = "2016-08-24 20:51:37"
# Error Detail cleared (handback version now in sync, no error)
$wsDeDe.Range("P2").Value = ""

# --- Column width adjustments (autofit-style widen for longer text) ---
# (ColumnWidth is expressed in character units and gets rounded to whole
# pixels by Excel when the file is saved, so the inputs below are chosen so
# the stored column width lands as close as possible to the target widths.)
$wsOverview.Range("E1").ColumnWidth = 29.15
$wsOverview.Range("F1").ColumnWidth = 29.15

$wsZhCn.Range("C1").ColumnWidth = 29.15
$wsZhCn.Range("P1").ColumnWidth = 12.83

$wsDeDe.Range("C1").ColumnWidth = 29.15
$wsDeDe.Range("P1").ColumnWidth = 12.83
